$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet (drop the stray " (2)" suffix) ---
$ws.Name = "repayment_20250916_20250916"

# --- Numeric (plain number) cell updates ---
$ws.Range("D2").Value = 4
$ws.Range("H2").Value = 1.409
$ws.Range("H3").Value = 477
$ws.Range("H4").Value = 550
$ws.Range("H5").Value = 470
$ws.Range("H6").Value = 410
$ws.Range("D7").Value = 4
$ws.Range("H7").Value = 152
$ws.Range("H8").Value = 233
$ws.Range("H9").Value = 639
$ws.Range("J9").Value = 1
$ws.Range("H10").Value = 419
$ws.Range("H11").Value = 700
$ws.Range("H12").Value = 442
$ws.Range("H13").Value = 216
$ws.Range("H14").Value = 786
$ws.Range("H15").Value = 1.179
$ws.Range("D16").Value = 8
$ws.Range("H16").Value = 448

# --- Text cell updates (values are stored as literal text strings,   ---
# --- same as the original cells, not auto-converted numbers/dates)   ---
$textCells = "E2","G2","E7","G7","K9","L9","E16","G16","K16"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "2,991,200.00"
$ws.Range("G2").Value = "1.89"
$ws.Range("E7").Value = "6,511,624.00"
$ws.Range("G7").Value = "3.72"
$ws.Range("K9").Value = "4.59"
$ws.Range("L9").Value = "7.14"
$ws.Range("E16").Value = "1,729,613.00"
$ws.Range("G16").Value = "1.22"
$ws.Range("K16").Value = "5.16"
